# Regenerate the "K" column (column G) values for the save_data sheet.
# These values are produced externally (K replaces the old "Strike#" calc,
# std/mean regenerated, s_vals recalculated) and simply written as literals
# into column G for each data row (rows 2-59), row 34 stays at 0 (unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, one per row starting at row 2 through row 59.
$newValues = @(
    1, 1, 1, 1, 2, 0, 0, 0, 1, 3,
    1, 0, 0, 0, 1, 3, 0, 1, 0, 0,
    0, 1, 1, 1, 1, 1, 3, 2, 0, 0,
    1, 1, 0, 1, 2, 1, 1, 2, 3, 0,
    0, 1, 1, 0, 1, 0, 2, 1, 0, 2,
    0, 0, 2, 1, 2, 3, 0, 2
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
